$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Nieuw-West " neighborhood label (with a trailing space) in A11 was a
# duplicate of "Nieuw-West" (no trailing space) used elsewhere in column A.
# Fix the typo so both rows share the same shared-string entry, which in turn
# drops the now-unused "Nieuw-West " string from the shared strings table.
$ws.Range("A11").Value = "Nieuw-West"

# Reflect the resulting selection/scroll position left behind by the edit.
$ws.Range("A11").Select()
